$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-06 16:35:21"
$wsOverview.Range("G3").Value = "2016-09-06 16:35:21"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-06 16:35:07"
$wsZhCn.Range("H3").Value = "2016-09-06 16:35:07"
$wsZhCn.Range("K2").Value = "2016-09-06 16:35:45"
$wsZhCn.Range("K3").Value = "2016-09-06 16:35:45"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-09-06 16:35:21"
$wsDeDe.Range("H3").Value = "2016-09-06 16:35:21"
$wsDeDe.Range("K2").Value = "2016-09-06 16:35:53"
$wsDeDe.Range("K3").Value = "2016-09-06 16:35:53"
